$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Test
$ws2 = $wb.Worksheets.Item(2)   # Configuration

# ---------------------------------------------------------------------------
# 1. Add the new "FindFlightData" sheet right after "Configuration"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $ws2)
$ws3.Name = "FindFlightData"

# ---------------------------------------------------------------------------
# 2. "Test" sheet - update Login row to RunMode=Yes and append a new
#    "Find Flight" data-driven test case
# ---------------------------------------------------------------------------
$ws1.Range("B2").Value = "Yes"

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Yes"
$ws1.Range("C4").Value = "Find Flight"
$ws1.Range("D4").Value = "tests.flight.FlightTest"
$ws1.Range("E4").Value = "P1,P2,P3"

# ---------------------------------------------------------------------------
# 3. "Configuration" sheet - insert a ThreadCount row before Environment
# ---------------------------------------------------------------------------
$ws2.Rows.Item(3).Insert()
$ws2.Range("A3").Value = "ThreadCount"
$ws2.Range("B3").Value = 1

# ---------------------------------------------------------------------------
# 4. "FindFlightData" sheet - column widths + header/data
# ---------------------------------------------------------------------------
$ws3.Columns.Item(2).ColumnWidth = 16.166666666666668
$ws3.Columns.Item(3).ColumnWidth = 15.736979166666666
$ws3.Columns.Item(4).ColumnWidth = 21.022135416666668
$ws3.Columns.Item(5).ColumnWidth = 11.877604166666666
$ws3.Columns.Item(6).ColumnWidth = 7.877604166666667
$ws3.Columns.Item(7).ColumnWidth = 21.592447916666668
$ws3.Columns.Item(8).ColumnWidth = 12.592447916666666
$ws3.Columns.Item(9).ColumnWidth = 10.022135416666666
$ws3.Columns.Item(10).ColumnWidth = 15.877604166666666
$ws3.Columns.Item(11).ColumnWidth = 24.022135416666668

$ws3.Range("A1").Value = "SrNo"
$ws3.Range("B1").Value = "Trip Type"
$ws3.Range("C1").Value = "No.Of Passengers"
$ws3.Range("D1").Value = "Source"
$ws3.Range("E1").Value = "Start Month"
$ws3.Range("F1").Value = "Start Day"
$ws3.Range("G1").Value = "Destination"
$ws3.Range("H1").Value = "Return Month"
$ws3.Range("I1").Value = "Return Day"
$ws3.Range("J1").Value = "Class"
$ws3.Range("K1").Value = "Airline Preference"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "Round Trip"
$ws3.Range("C2").Value = 1
$ws3.Range("D2").Value = "Frankfurt"
$ws3.Range("E2").Value = "February"
$ws3.Range("F2").Value = 9
$ws3.Range("G2").Value = "London"
$ws3.Range("H2").Value = "February"
$ws3.Range("I2").Value = 19
$ws3.Range("J2").Value = "Economy"

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "One Way"
$ws3.Range("C3").Value = 2
$ws3.Range("D3").Value = "Paris"
$ws3.Range("E3").Value = "June"
$ws3.Range("F3").Value = 7
$ws3.Range("G3").Value = "Seattle"
$ws3.Range("J3").Value = "Business"
$ws3.Range("K3").Value = "Blue Skies Airlines"

$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = "Round Trip"
$ws3.Range("C4").Value = 3
$ws3.Range("D4").Value = "Seattle"
$ws3.Range("E4").Value = "March"
$ws3.Range("F4").Value = 14
$ws3.Range("G4").Value = "Acapulco"
$ws3.Range("H4").Value = "March"
$ws3.Range("I4").Value = 20
$ws3.Range("J4").Value = "First"
$ws3.Range("K4").Value = "Unified Airlines"

$ws3.Range("A5").Value = 4
$ws3.Range("B5").Value = "One Way"
$ws3.Range("C5").Value = 4
$ws3.Range("D5").Value = "Sydney"
$ws3.Range("E5").Value = "October"
$ws3.Range("F5").Value = 30
$ws3.Range("G5").Value = "Portland"
$ws3.Range("J5").Value = "Economy"
$ws3.Range("K5").Value = "Pangea Airlines"

# ---------------------------------------------------------------------------
# 5. Formatting - bold+border header rows, bordered data rows, on every sheet
# ---------------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2, $ws3)) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    $lastCol = $used.Columns.Count

    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
    $headerRange.Font.Bold = $true
    $headerRange.Borders.LineStyle = 1

    if ($lastRow -gt 1) {
        $dataRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, $lastCol))
        $dataRange.Borders.LineStyle = 1
    }
}

# ---------------------------------------------------------------------------
# 6. Selections / active tab - Configuration becomes the active sheet
# ---------------------------------------------------------------------------
$ws1.Range("H9").Select()
$ws3.Range("I9").Select()

$ws2.Activate()
$ws2.Range("B3").Select()
